$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "Impact" analysis columns (C, E, G) that are no longer used
$ws.Range("C1").ClearContents()
$ws.Range("E1").ClearContents()
$ws.Range("G1").ClearContents()

$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()

$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()

$ws.Range("C7").ClearContents()
$ws.Range("E7").ClearContents()

$ws.Range("C9").ClearContents()
$ws.Range("E9").ClearContents()

$ws.Range("C11").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C17").ClearContents()

# Row 3 loses its explicit wrap-driven height once the wrapped text is gone
$ws.Rows.Item(3).AutoFit()

# Move the selection
$ws.Range("C9").Select()
